$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "게스트가 필수 입력 정보 (ID, 비밀`n번호, 전화번호, 결제 수단, 선호 자전거 유형(일반/전기) 등)을 입력하여 회원 계정을 생성하는 기능"
$ws.Range("C8").Value = "회원 가입"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "회원 정보를 회원 리스트에서 삭제하는 기능"
$ws.Range("C9").Value = "회원 탈퇴"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "회원이 ID와 비밀번호를 입력하고 시스템에 접속한다."
$ws.Range("C10").Value = "로그인하기"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "로그아웃으로 시스템 접속을 종료한다."
$ws.Range("C11").Value = "로그아웃하기"
